$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1065.3103
$ws.Range("J17").Value = 1100.1786
$ws.Range("L17").Value = 3300.5358
$ws.Range("N17").Value = -3636.5358
$ws.Range("H74").Value = 3326.25
$ws.Range("I74").Value = 3326.25
$ws.Range("K74").Value = 3326.25
$ws.Range("M74").Value = -2390.25
$ws.Range("H77").Value = 3326.25
$ws.Range("I77").Value = 3326.25
$ws.Range("K77").Value = 16631.25
$ws.Range("M77").Value = -11951.25
$ws.Range("H132").Value = 10104045
$ws.Range("I132").Value = 14496083
$ws.Range("J132").Value = 2356.6
$ws.Range("K132").Value = 43488249
$ws.Range("L132").Value = 7069.799999999999
$ws.Range("M132").Value = -43485719
$ws.Range("N132").Value = -12129.8
$ws.Range("H137").Value = 2041.4062
$ws.Range("I137").Value = 1134.0667
$ws.Range("J137").Value = 2842
$ws.Range("K137").Value = 3402.2001
$ws.Range("L137").Value = 8526
$ws.Range("M137").Value = -852.2001
$ws.Range("N137").Value = -13626
$ws.Range("H138").Value = 1635.2347
$ws.Range("I138").Value = 660.44446
$ws.Range("J138").Value = 1854.5625
$ws.Range("K138").Value = 1981.33338
$ws.Range("L138").Value = 5563.6875
$ws.Range("M138").Value = 3158.66662
$ws.Range("N138").Value = -15843.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2829.849
$ws.Range("I32").Value = 3067.805
$ws.Range("J32").Value = 2016.8334
$ws.Range("K32").Value = 3067.805
$ws.Range("L32").Value = 2016.8334
$ws.Range("M32").Value = -2780.805
$ws.Range("N32").Value = -2590.8334
$ws.Range("H61").Value = 1633.4
$ws.Range("I61").Value = 1633.4
$ws.Range("K61").Value = 1633.4
$ws.Range("M61").Value = -1421.4
$ws.Range("H63").Value = 125001940
$ws.Range("I63").Value = 2300
$ws.Range("K63").Value = 2300
$ws.Range("M63").Value = -1614
$ws.Range("H66").Value = 125001940
$ws.Range("I66").Value = 2300
$ws.Range("K66").Value = 11500
$ws.Range("M66").Value = -8068
$ws.Range("H74").Value = 902.7059
$ws.Range("I74").Value = 770.1429000000001
$ws.Range("J74").Value = 1521.3334
$ws.Range("K74").Value = 770.1429000000001
$ws.Range("L74").Value = 1521.3334
$ws.Range("M74").Value = 103.8570999999999
$ws.Range("N74").Value = -3269.3334
$ws.Range("H77").Value = 902.7059
$ws.Range("I77").Value = 770.1429000000001
$ws.Range("J77").Value = 1521.3334
$ws.Range("K77").Value = 3850.7145
$ws.Range("L77").Value = 7606.666999999999
$ws.Range("M77").Value = 517.2855
$ws.Range("N77").Value = -16342.667
$ws.Range("H97").Value = 493.0435
$ws.Range("I97").Value = 422.57895
$ws.Range("J97").Value = 827.75
$ws.Range("K97").Value = 422.57895
$ws.Range("L97").Value = 827.75
$ws.Range("M97").Value = 73.42104999999998
$ws.Range("N97").Value = -1819.75
$ws.Range("H136").Value = 1633.4
$ws.Range("I136").Value = 1633.4
$ws.Range("K136").Value = 4900.200000000001
$ws.Range("M136").Value = -2350.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 150.71428
$ws.Range("I22").Value = 142.5
$ws.Range("K22").Value = 142.5
$ws.Range("M22").Value = 30.5
$ws.Range("H134").Value = 16794.715
$ws.Range("I134").Value = 11739.728
$ws.Range("J134").Value = 35329.668
$ws.Range("K134").Value = 35219.18399999999
$ws.Range("L134").Value = 105989.004
$ws.Range("M134").Value = -32684.18399999999
$ws.Range("N134").Value = -111059.004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 10753892
$ws.Range("I134").Value = 12821545
$ws.Range("K134").Value = 38464635
$ws.Range("M134").Value = -38462100
$ws.Range("H140").Value = 39000
$ws.Range("J140").Value = 39000
$ws.Range("L140").Value = 39000
$ws.Range("N140").Value = -49360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2083.3333
$ws.Range("H93").Value = 4996.7915
$ws.Range("J93").Value = 4996.7915
$ws.Range("L93").Value = 14990.3745
$ws.Range("N93").Value = -18734.3745
$ws.Range("H114").Value = 727.56525
$ws.Range("I114").Value = 420.8889
$ws.Range("J114").Value = 924.7143
$ws.Range("K114").Value = 1262.6667
$ws.Range("L114").Value = 2774.1429
$ws.Range("M114").Value = 1991.3333
$ws.Range("N114").Value = -9282.142899999999
$ws.Range("H117").Value = 918
$ws.Range("J117").Value = 1174
$ws.Range("L117").Value = 3522
$ws.Range("N117").Value = -10406
$ws.Range("H122").Value = 671.125
$ws.Range("I122").Value = 562
$ws.Range("J122").Value = 853
$ws.Range("K122").Value = 5058
$ws.Range("L122").Value = 7677
$ws.Range("M122").Value = -2608
$ws.Range("N122").Value = -12577
$ws.Range("H129").Value = 14368965
$ws.Range("J129").Value = 4167991
$ws.Range("L129").Value = 12503973
$ws.Range("N129").Value = -12513973
$ws.Range("H131").Value = 16667966
$ws.Range("I131").Value = 100000680
$ws.Range("J131").Value = 1423.2
$ws.Range("K131").Value = 300002040
$ws.Range("L131").Value = 4269.6
$ws.Range("M131").Value = -299997000
$ws.Range("N131").Value = -14349.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 23688958
$ws.Range("J70").Value = 20005630
$ws.Range("L70").Value = 20005630
$ws.Range("N70").Value = -20006170
$ws.Range("H73").Value = 23688958
$ws.Range("J73").Value = 20005630
$ws.Range("L73").Value = 20005630
$ws.Range("N73").Value = -20007502
$ws.Range("H102").Value = 2637.4
$ws.Range("I102").Value = 3408.7144
$ws.Range("J102").Value = 1962.5
$ws.Range("K102").Value = 3408.7144
$ws.Range("L102").Value = 1962.5
$ws.Range("M102").Value = -1786.7144
$ws.Range("N102").Value = -5206.5
$ws.Range("H132").Value = 3125.4119
$ws.Range("I132").Value = 2759.5715
$ws.Range("K132").Value = 8278.7145
$ws.Range("M132").Value = -5748.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1612.75
$ws.Range("I68").Value = 1240.8
$ws.Range("K68").Value = 1240.8
$ws.Range("M68").Value = -491.8
$ws.Range("H71").Value = 1612.75
$ws.Range("I71").Value = 1240.8
$ws.Range("K71").Value = 6204
$ws.Range("M71").Value = -2460
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H136").Value = 9584.916999999999
$ws.Range("I136").Value = 13625.625
$ws.Range("K136").Value = 40876.875
$ws.Range("M136").Value = -38326.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 45460164
$ws.Range("I62").Value = 55559776
$ws.Range("J62").Value = 11899.5
$ws.Range("K62").Value = 55559776
$ws.Range("L62").Value = 11899.5
$ws.Range("M62").Value = -55559152
$ws.Range("N62").Value = -13147.5
$ws.Range("H65").Value = 45460164
$ws.Range("I65").Value = 55559776
$ws.Range("J65").Value = 11899.5
$ws.Range("K65").Value = 277798880
$ws.Range("L65").Value = 59497.5
$ws.Range("M65").Value = -277795760
$ws.Range("N65").Value = -65737.5
$ws.Range("H107").Value = 531.4211
$ws.Range("I107").Value = 462.125
$ws.Range("K107").Value = 1386.375
$ws.Range("M107").Value = 533.625
$ws.Range("H122").Value = 37153372
$ws.Range("I122").Value = 43345268
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 130035804
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -130033354
$ws.Range("N122").Value = -10900
$ws.Range("H136").Value = 1030
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
